# Update column G ("K") values for rows 2-14 on the active worksheet.
# These values replace the previous "Strike#" based values with the
# regenerated "K" values described in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 1
    3  = 0
    4  = 2
    5  = 3
    6  = 0
    7  = 1
    8  = 0
    9  = 3
    10 = 0
    11 = 1
    12 = 2
    13 = 0
    14 = 2
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
